$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Package info": bump version/build metadata and add a "Name" row.
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Package info")

$wsInfo.Cells.Item(4,2).Value = "1.1.0"
$wsInfo.Cells.Item(5,2).Value = "2.33.9"
$wsInfo.Cells.Item(6,1).Value = "DHIS2 build"
$wsInfo.Cells.Item(6,2).Value = "58094d2"
$wsInfo.Cells.Item(7,1).Value = "Last updated"
$wsInfo.Cells.Item(7,2).Value = "20210520T090044"

# New row 8 - copy formatting from row 6 (same banding) then set the values.
$wsInfo.Range("A6:B6").Copy()
$wsInfo.Range("A8:B8").PasteSpecial(-4122)
$wsInfo.Cells.Item(8,1).Value = "Name"
$wsInfo.Cells.Item(8,2).Value = "MAL-HIST_CUSTOM_V1.1.0_2.33.9-en"

# Narrow column B.
$wsInfo.Columns.Item(2).ColumnWidth = 33.83

# ---------------------------------------------------------------------------
# Sheet "dataElements": rows 2-7 get reshuffled into a new order.
# ---------------------------------------------------------------------------
$wsDE = $wb.Worksheets.Item("dataElements")

$wsDE.Cells.Item(2,1).Value = "MAL - Plasmodium falciparum (Mic+RDT)"
$wsDE.Cells.Item(2,2).Value = "P. falciparum (micr+RDT)"
$wsDE.Cells.Item(2,3).Value = "MAL_PF_MICR_RDT"
$wsDE.Cells.Item(2,4).Value = "Cases confirmed as P.falciparum positive with microscopy and/or RDT"
$wsDE.Cells.Item(2,5).Value = "IvYR8mc6prX"
$wsDE.Cells.Item(2,6).Value = "2019-10-20"
$wsDE.Cells.Item(2,7).Value = "IIU1O0Z4l49"

$wsDE.Cells.Item(3,1).Value = "MAL - Mixed/Other malaria species (Mic+RDT)"
$wsDE.Cells.Item(3,2).Value = "Other species (micr+RDT)"
$wsDE.Cells.Item(3,3).Value = "MAL_MIX_OTHER_SPECIES_MICR_RDT"
$wsDE.Cells.Item(3,4).Value = "Cases confirmed as P.malariae or P. ovale or P. knowlesi with microscopy and/or RDT"
$wsDE.Cells.Item(3,5).Value = "IvYR8mc6prX"
$wsDE.Cells.Item(3,6).Value = "2019-10-20"
$wsDE.Cells.Item(3,7).Value = "JkOyLRb3dpX"

$wsDE.Cells.Item(4,1).Value = "MAL - Plasmodium vivax (Mic+RDT)"
$wsDE.Cells.Item(4,2).Value = "P. vivax (micr+RDT)"
$wsDE.Cells.Item(4,3).Value = "MAL_PV_MICR_RDT"
$wsDE.Cells.Item(4,4).Value = "Cases confirmed as P.vivax positive with microscopy"
$wsDE.Cells.Item(4,5).Value = "IvYR8mc6prX"
$wsDE.Cells.Item(4,6).Value = "2019-10-20"
$wsDE.Cells.Item(4,7).Value = "pUC8tgzn0lV"

$wsDE.Cells.Item(5,1).Value = "MAL - Mixed malaria species  (Mic+RDT)"
$wsDE.Cells.Item(5,2).Value = "Mixed (micr+RDT)"
$wsDE.Cells.Item(5,3).Value = "MAL_MIX_SPECIES_MICR_RDT"
$wsDE.Cells.Item(5,4).Value = "Cases confirmed as mixed infection with P.falciparum and P.vivax with microscopy and/or RDT"
$wsDE.Cells.Item(5,5).Value = "IvYR8mc6prX"
$wsDE.Cells.Item(5,6).Value = "2019-10-20"
$wsDE.Cells.Item(5,7).Value = "TNTW2ruEVEu"

$wsDE.Cells.Item(6,1).Value = "MAL - Malaria tested cases (Mic+RDT)"
$wsDE.Cells.Item(6,2).Value = "Tested (micr+RDT)"
$wsDE.Cells.Item(6,3).Value = "MAL_TEST_CASES_MICR_RDT"
$wsDE.Cells.Item(6,4).Value = "Suspected cases tested with both microscopy and/or RDT"
$wsDE.Cells.Item(6,5).Value = "IvYR8mc6prX"
$wsDE.Cells.Item(6,6).Value = "2019-10-20"
$wsDE.Cells.Item(6,7).Value = "tuOTgWfDO6m"

$wsDE.Cells.Item(7,1).Value = "MAL - Malaria confirmed cases (Mic+RDT)"
$wsDE.Cells.Item(7,2).Value = "Positive (micr+RDT)"
$wsDE.Cells.Item(7,3).Value = "MAL_CONFI_CASES_MICR_RDT"
$wsDE.Cells.Item(7,4).Value = "Cases confirmed as positive with microscopy and/or RDT"
$wsDE.Cells.Item(7,5).Value = "IvYR8mc6prX"
$wsDE.Cells.Item(7,6).Value = "2019-10-20"
$wsDE.Cells.Item(7,7).Value = "X0luAFiy268"

# ---------------------------------------------------------------------------
# Sheet "dataElementGroups": column B follows the same reshuffle (column A
# stays "Malaria old records only" throughout).
# ---------------------------------------------------------------------------
$wsDEG = $wb.Worksheets.Item("dataElementGroups")

$wsDEG.Cells.Item(2,2).Value = "MAL - Plasmodium falciparum (Mic+RDT)"
$wsDEG.Cells.Item(3,2).Value = "MAL - Mixed/Other malaria species (Mic+RDT)"
$wsDEG.Cells.Item(4,2).Value = "MAL - Plasmodium vivax (Mic+RDT)"
$wsDEG.Cells.Item(5,2).Value = "MAL - Mixed malaria species  (Mic+RDT)"
$wsDEG.Cells.Item(6,2).Value = "MAL - Malaria tested cases (Mic+RDT)"
$wsDEG.Cells.Item(7,2).Value = "MAL - Malaria confirmed cases (Mic+RDT)"

# ---------------------------------------------------------------------------
# Sheet "userGroups": rows 2 and 4 swap (Name/UID), all rows get the new
# "Last updated" date.
# ---------------------------------------------------------------------------
$wsUG = $wb.Worksheets.Item("userGroups")

$wsUG.Cells.Item(2,1).Value = "Malaria data capture"
$wsUG.Cells.Item(2,2).Value = "2021-05-20"
$wsUG.Cells.Item(2,3).Value = "fRSrUJ6SMGH"

$wsUG.Cells.Item(3,1).Value = "Malaria admin"
$wsUG.Cells.Item(3,2).Value = "2021-05-20"
$wsUG.Cells.Item(3,3).Value = "suMb19wGXPR"

$wsUG.Cells.Item(4,1).Value = "Malaria access"
$wsUG.Cells.Item(4,2).Value = "2021-05-20"
$wsUG.Cells.Item(4,3).Value = "ZXEVDM9XRea"
